$d = $word.ActiveDocument

# The sentence originally reads:
#   "...uma contribuição do Fábio Orsi."
# and the runs are split as:
#   "contribuição do Fábio" | " Orsi"
# After the edit it must read:
#   "...uma contribuição do Fábio."
# with the runs split as:
#   "contribuição do Fábi" | "o"

# Step 1: remove the trailing " Orsi" (keeping the final period).
$delRange = $d.Content
$delRange.Find.Execute(" Orsi", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

# Step 2: force the remaining "Fábio" run to split into "Fábi" + "o" by
# nudging the formatting of the last letter and reverting it, which makes
# the run boundary persist without altering the visible formatting.
$searchRange = $d.Content
$searchRange.Find.Execute("contribuição do Fábio", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$lastLetter = $d.Range($searchRange.End - 1, $searchRange.End)
$lastLetter.Bold = 1
$lastLetter.Bold = 0
